$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update "想去人数" (column F) for rows 4-10
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 220
$ws1.Range("F5").Value = 2656
$ws1.Range("F6").Value = 1884
$ws1.Range("F7").Value = 364
$ws1.Range("F8").Value = 113
$ws1.Range("F9").Value = 938
$ws1.Range("F10").Value = 181

# Sheet "全部类型" (All types) - update "想去人数" (column F) for rows 4-11 (row 8 unchanged)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 220
$ws4.Range("F5").Value = 2656
$ws4.Range("F6").Value = 1884
$ws4.Range("F7").Value = 364
$ws4.Range("F9").Value = 113
$ws4.Range("F10").Value = 938
$ws4.Range("F11").Value = 181
